$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: LDAP paragraph - the run split "...a la q" / "ue pueden..."
# becomes a single run "...a la que pueden realizarse consultas."
# (plain text is identical across the run boundary, so a find/replace
# with the same text simply causes Word to re-flow/merge the runs)
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "a la que pueden realizarse consultas.", $true, $false, $false, $false, $false,
    $true, 1, $false, "a la que pueden realizarse consultas.", 2) | Out-Null

# -----------------------------------------------------------------
# Change 3: SMTP paragraph - the run split "...y es u" / "n estándar..."
# becomes a single run "...y es un estándar oficial de Internet."
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "es un estándar oficial de Internet.", $true, $false, $false, $false, $false,
    $true, 1, $false, "es un estándar oficial de Internet.", 2) | Out-Null

# -----------------------------------------------------------------
# Change 2: POP paragraph - merge "...de nivel" / " de aplicación en el
# Modelo OSI" into a single run.
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "Es un protocolo de nivel de aplicación en el Modelo OSI", $true, $false, $false, $false, $false,
    $true, 1, $false, "Es un protocolo de nivel de aplicación en el Modelo OSI", 2) | Out-Null

# -----------------------------------------------------------------
# Change 2 (cont.): the paragraph that used to hold the _GoBack
# bookmark (right after the POP section, before the SMTP heading)
# becomes a plain empty paragraph - the bookmark moves to the very
# end of the document.
# -----------------------------------------------------------------
$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $xml = $p.Range.WordOpenXML
    if ($xml -like "*_GoBack*") {
        $body = '<w:p/>'
        $p.Range.InsertXML($pkgHeader + "<w:body>$body</w:body>" + $pkgFooter)
        $found = $true
        break
    }
}

# -----------------------------------------------------------------
# Change 4: append the new "Registros DNS" section before the final
# (empty) paragraph of the document, then turn that final paragraph
# into the new home of the _GoBack bookmark.
# -----------------------------------------------------------------
$content = '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Registros DNS</w:t></w:r></w:p>' +
    '<w:p><w:r><w:t xml:space="preserve">Los registros DNS son archivos de mapeo o sistemas que le indican a un servidor DNS a qué dirección IP está asociado un dominio particular. También </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>le</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> indican a los servidores DNS cómo manejar las solicitudes que se envían a cada nombre de dominio.</w:t></w:r></w:p>' +
    '<w:p/><w:p/><w:p/><w:p/><w:p/>' +
    '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>¿Qué tipo de registros DNS existen y para qué sirven?</w:t></w:r></w:p>' +
    '<w:p><w:r><w:t>El lugar donde se configuran las entradas DNS para cada dominio son los servidores de nombres. Los diferentes tipos de entradas de registro son:</w:t></w:r></w:p>' +
    '<w:p/>' +
    '<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Registro A</w:t></w:r><w:r><w:t>: Este registro se utiliza para convertir nombres de host en direcciones IP.</w:t></w:r></w:p>' +
    '<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Registro CNAME</w:t></w:r><w:r><w:t>: Se utiliza para crear nombres de host adicionales (alias), y para crear diferentes servicios bajo una misma dirección IP.</w:t></w:r></w:p>' +
    '<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Registro NS</w:t></w:r><w:r><w:t>: indica los servidores de DNS autorizados para el dominio, es decir, a quién hay que preguntar para saber acerca de los registros de midominio.info.</w:t></w:r></w:p>' +
    '<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Registro MX</w:t></w:r><w:r><w:t>: Se utiliza para asociar un nombre de dominio a una lista de servidores de correo para la recepción de emails. Nos interesa si queremos realizar redirecciones de nuestro correo o utilizar nuestro correo electrónico con otro proveedor.</w:t></w:r></w:p>' +
    '<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Registro SPF</w:t></w:r><w:r><w:t>: define qué servidores están autorizados para enviar correo electrónico con nuestro dominio.</w:t></w:r></w:p>' +
    '<w:p><w:r><w:t>Configurando adecuadamente estos registros podemos exprimir al máximo todas las funcionalidades que poseen las DNS de nuestro dominio.</w:t></w:r></w:p>' +
    '<w:p/>'

$paras = $d.Paragraphs
$lastPara = $paras.Item($paras.Count)
$lastRange = $lastPara.Range
$lastRange.Collapse(1)
$lastRange.InsertXML($pkgHeader + "<w:body>$content</w:body>" + $pkgFooter)

# Now the true final (still-empty) paragraph gets the relocated bookmark.
$paras2 = $d.Paragraphs
$finalPara = $paras2.Item($paras2.Count)
$finalRange = $finalPara.Range
$finalRange.Collapse(1)
$finalRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')

Write-Output "done"
